$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(0, 4, 1, 2, 4, 2, 4, 8, 3, 1, 2, 1, 4, 4, 1, 1)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Range("G$row").Value = $values[$i]
}
